# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the handback
# file for 4b732e66-d30b-4fde-b4a3-e04f9979acbe.md is out of date and the
# item is being sent back out ("Ready for handoff") for both zh-cn and de-de.

$wb = $excel.ActiveWorkbook

$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfa0804bf6c3663c9148544a1b46cd4a31867359/e2e/4b732e66-d30b-4fde-b4a3-e04f9979acbe.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e49197dfb9973bc8f5fcec1ea7db53423aa5a13/e2e/4b732e66-d30b-4fde-b4a3-e04f9979acbe.md.'

# --- Overview sheet: row 3 is the 4b732e66-...-md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 04:47:46"

# --- zh-cn sheet: row 3 is the 4b732e66-...-md file ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("H3").Value = "2016-08-28 04:47:42"
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is the 4b732e66-...-md file ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("H3").Value = "2016-08-28 04:47:46"
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
